$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 2 (C2:F2) previously held empty, highlighted cells. Fill them with the
# Android-locator strings used by the new reset-pwd test case / control updates.
$ws.Range("C2").Value = "//android.widget.Spinner[@hint='VON:']"
$ws.Range("D2").Value = "//android.widget.Spinner[@hint='BIS:']"
$ws.Range("E2").Value = "//android.widget.EditText[@hint='KOMMENTAR (NUR FÜR DICH SICHTBAR):']"
$ws.Range("F2").Value = '//android.widget.Button[@text="Speichern"]'

# Selection moved from D5 to F6 as part of the editing session.
$ws.Range("F6").Select()
